$wb = $excel.ActiveWorkbook

# 1) Rename sheet "DAY 5 (08-04-2022) (2)" -> "DAY 5 (08-04-2022)"
$ws5 = $wb.Worksheets.Item("DAY 5 (08-04-2022) (2)")
$ws5.Name = "DAY 5 (08-04-2022)"

# 2) Edit cells on "Day 7 (11-04-2022)" sheet
$ws = $wb.Worksheets.Item("Day 7 (11-04-2022)")

# D9: "...Data Modelling..." -> "...Discussed Data Modelling..."
$ws.Range("D9").Value = "1 Hr 30 Mins: Brainstorming                                             1 Hr 30 mins: Discussed Data Modelling                                           1 Hr :Entered attribute values(Employee,TAC)"

# D15: adjust durations in Sheik Fareeth's entry
$ws.Range("D15").Value = "1 hr 30 mins - Brain Stromming`n30 Mins - Adding New slides to TAC ( Upcoming drives, notifications, scheduling drives - 5 slides )`n1 hr 30 mins - Building Data model In draw.io ( 7 Entities )`n30 Mins - Re refining Interviewers scheduled drive cancellation"

# D16: replace Vinoth's entry with new text, and bump hours F16 2.5 -> 4
$ws.Range("D16").Value = "1Hr 30 mins : Brainstroming`n1 hr : Discussed Data modeling`n1 hr 30 mins : Entered attributes for data model(Pool managing, invites, interviews)                                                           "
$ws.Range("F16").Value = 4

# 3) Update the sheet view: scroll position + active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G16").Select()

